$wb = $excel.ActiveWorkbook

# --- ALC (diff hunk @ line 3202) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 280.78946
$ws.Range("I53").Value = 237.53847
$ws.Range("J53").Value = 374.5
$ws.Range("K53").Value = 237.53847
$ws.Range("L53").Value = 374.5
$ws.Range("M53").Value = 399.46153
$ws.Range("N53").Value = -1648.5

# --- ALC (diff hunk @ line 5911) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 16130859
$ws.Range("I107").Value = 10418166
$ws.Range("K107").Value = 10418166
$ws.Range("M107").Value = -10416246

# --- ALC (diff hunk @ line 7145) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2960.625
$ws.Range("J132").Value = 7444.727
$ws.Range("L132").Value = 22334.181
$ws.Range("N132").Value = -27394.181

# --- ALC (diff hunk @ line 7295) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 3202.0667
$ws.Range("I135").Value = 1944.4
$ws.Range("J135").Value = 5717.4
$ws.Range("K135").Value = 17499.6
$ws.Range("L135").Value = 51456.6
$ws.Range("M135").Value = -14964.6
$ws.Range("N135").Value = -56526.6

# --- ALC (diff hunk @ line 7396) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 45921.04
$ws.Range("I137").Value = 54526
$ws.Range("K137").Value = 163578
$ws.Range("M137").Value = -161028

# --- ALC (diff hunk @ line 7448) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2674.04
$ws.Range("I138").Value = 1206.0526
$ws.Range("J138").Value = 3573.7742
$ws.Range("K138").Value = 3618.1578
$ws.Range("L138").Value = 10721.3226
$ws.Range("M138").Value = 1521.8422
$ws.Range("N138").Value = -21001.3226

# --- ARM (diff hunk @ line 9199) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7475.176
$ws.Range("I32").Value = 7475.176
$ws.Range("K32").Value = 7475.176
$ws.Range("M32").Value = -7188.176

# --- ARM (diff hunk @ line 10608) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7859.08
$ws.Range("I61").Value = 4199.5386
$ws.Range("K61").Value = 4199.5386
$ws.Range("M61").Value = -3987.5386

# --- ARM (diff hunk @ line 11230) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 53135.332
$ws.Range("I74").Value = 53135.332
$ws.Range("K74").Value = 53135.332
$ws.Range("M74").Value = -52261.332

# --- ARM (diff hunk @ line 11377) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 53135.332
$ws.Range("I77").Value = 53135.332
$ws.Range("K77").Value = 265676.66
$ws.Range("M77").Value = -261308.66

# --- ARM (diff hunk @ line 14045) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5524.2104
$ws.Range("I132").Value = 4338.6875
$ws.Range("K132").Value = 13016.0625
$ws.Range("M132").Value = -10486.0625

# --- ARM (diff hunk @ line 14238) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7859.08
$ws.Range("I136").Value = 4199.5386
$ws.Range("K136").Value = 12598.6158
$ws.Range("M136").Value = -10048.6158

# --- BSM (diff hunk @ line 21013) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3120.611
$ws.Range("I134").Value = 2585.0715
$ws.Range("J134").Value = 4995
$ws.Range("K134").Value = 7755.2145
$ws.Range("L134").Value = 14985
$ws.Range("M134").Value = -5220.2145
$ws.Range("N134").Value = -20055

# --- CRP (diff hunk @ line 22206) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3545.111
$ws.Range("I16").Value = 3501
$ws.Range("K16").Value = 3501
$ws.Range("M16").Value = -3214

# --- CRP (diff hunk @ line 22920) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 338520.47
$ws.Range("I31").Value = 478366.62
$ws.Range("J31").Value = 71541.45
$ws.Range("K31").Value = 478366.62
$ws.Range("L31").Value = 71541.45
$ws.Range("M31").Value = -478071.62
$ws.Range("N31").Value = -72131.45

# --- CRP (diff hunk @ line 23070) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 338520.47
$ws.Range("I34").Value = 478366.62
$ws.Range("J34").Value = 71541.45
$ws.Range("K34").Value = 478366.62
$ws.Range("L34").Value = 71541.45
$ws.Range("M34").Value = -478164.62
$ws.Range("N34").Value = -71945.45

# --- CRP (diff hunk @ line 24213) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1938
$ws.Range("I58").Value = 1983.25
$ws.Range("J58").Value = 1214
$ws.Range("K58").Value = 1983.25
$ws.Range("L58").Value = 1214
$ws.Range("M58").Value = -1780.25
$ws.Range("N58").Value = -1620

# --- CRP (diff hunk @ line 26620) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 3910.32
$ws.Range("I107").Value = 845.5454999999999
$ws.Range("K107").Value = 845.5454999999999
$ws.Range("M107").Value = 1074.4545

# --- CRP (diff hunk @ line 26914) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 3545.111
$ws.Range("I113").Value = 3501
$ws.Range("K113").Value = 3501
$ws.Range("M113").Value = -1331

# --- CRP (diff hunk @ line 27818) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4659
$ws.Range("I132").Value = 2915.1155
$ws.Range("J132").Value = 50000
$ws.Range("K132").Value = 8745.3465
$ws.Range("L132").Value = 150000
$ws.Range("M132").Value = -6215.3465
$ws.Range("N132").Value = -155060

# --- CRP (diff hunk @ line 27919) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6085.4814
$ws.Range("I134").Value = 6558.9546
$ws.Range("J134").Value = 4002.2
$ws.Range("K134").Value = 19676.8638
$ws.Range("L134").Value = 12006.6
$ws.Range("M134").Value = -17141.8638
$ws.Range("N134").Value = -17076.6

# --- CRP (diff hunk @ line 28020) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1938
$ws.Range("I136").Value = 1983.25
$ws.Range("J136").Value = 1214
$ws.Range("K136").Value = 5949.75
$ws.Range("L136").Value = 3642
$ws.Range("M136").Value = -3399.75
$ws.Range("N136").Value = -8742

# --- CUL (diff hunk @ line 28506) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 20042928
$ws.Range("I4").Value = 20042928
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 60128784
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -60128672
$ws.Range("N4").ClearContents()

# --- CUL (diff hunk @ line 32682) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 2200
$ws.Range("J88").Value = 2200
$ws.Range("L88").Value = 6600
$ws.Range("N88").Value = -7456

# --- CUL (diff hunk @ line 32832) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 2200
$ws.Range("J91").Value = 2200
$ws.Range("L91").Value = 6600
$ws.Range("N91").Value = -9564

# --- CUL (diff hunk @ line 34657) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 241333
$ws.Range("I128").Value = 241333
$ws.Range("K128").Value = 723999
$ws.Range("M128").Value = -719019

# --- GSM (diff hunk @ line 38790) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8581.826999999999
$ws.Range("I70").Value = 8632.823
$ws.Range("J70").Value = 8509.583000000001
$ws.Range("K70").Value = 8632.823
$ws.Range("L70").Value = 8509.583000000001
$ws.Range("M70").Value = -8362.823
$ws.Range("N70").Value = -9049.583000000001

# --- GSM (diff hunk @ line 38934) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8581.826999999999
$ws.Range("I73").Value = 8632.823
$ws.Range("J73").Value = 8509.583000000001
$ws.Range("K73").Value = 8632.823
$ws.Range("L73").Value = 8509.583000000001
$ws.Range("M73").Value = -7696.823
$ws.Range("N73").Value = -10381.583

# --- GSM (diff hunk @ line 39274) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3956
$ws.Range("I80").Value = 3871.5
$ws.Range("K80").Value = 3871.5
$ws.Range("M80").Value = -2873.5

# --- GSM (diff hunk @ line 39418) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3956
$ws.Range("I83").Value = 3871.5
$ws.Range("K83").Value = 19357.5
$ws.Range("M83").Value = -14365.5

# --- GSM (diff hunk @ line 40588) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 13910.75
$ws.Range("I107").Value = 17047.834
$ws.Range("J107").Value = 4499.5
$ws.Range("K107").Value = 17047.834
$ws.Range("L107").Value = 4499.5
$ws.Range("M107").Value = -15127.834
$ws.Range("N107").Value = -8339.5

# --- GSM (diff hunk @ line 41777) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 34835.234
$ws.Range("I132").Value = 39118.25
$ws.Range("K132").Value = 117354.75
$ws.Range("M132").Value = -114824.75

# --- LTW (diff hunk @ line 44964) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 650.5625
$ws.Range("I55").Value = 634.25
$ws.Range("J55").Value = 699.5
$ws.Range("K55").Value = 634.25
$ws.Range("L55").Value = 699.5
$ws.Range("M55").Value = -461.25
$ws.Range("N55").Value = -1045.5

# --- LTW (diff hunk @ line 48677) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8938.736999999999
$ws.Range("I132").Value = 10066.723
$ws.Range("J132").Value = 7923.55
$ws.Range("K132").Value = 30200.169
$ws.Range("L132").Value = 23770.65
$ws.Range("M132").Value = -27670.169
$ws.Range("N132").Value = -28830.65

# --- WVR (diff hunk @ line 55770) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 198065.72
$ws.Range("I136").Value = 228751.86
$ws.Range("J136").Value = 5181.4287
$ws.Range("K136").Value = 686255.58
$ws.Range("L136").Value = 15544.2861
$ws.Range("M136").Value = -20644.2861
